# Canopy cover generate module implemented:
# - Column G ("geometry") removed entirely (header + per-row MULTIPOLYGON values)
# - Column F (previously all 0, under header "20240109") now holds computed
#   canopy-cover percentage values per plot row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    36.15187382156221,
    43.18512393755113,
    14.17988118938711,
    14.46408889828333,
    6.229167795332166,
    12.04594757992816,
    23.8212848256756,
    38.26068594584124,
    21.43054609577148,
    36.85519139964804,
    42.17963890126777,
    33.0578973668903,
    46.37618400860777,
    36.20023636075,
    49.49106929131641,
    34.55137010425472,
    22.92929907279351,
    31.02294172764158,
    35.08699557419547,
    52.43948124387937,
    49.38161259357787,
    31.02058455188471,
    45.64289133893747,
    34.7198460975232,
    35.36163876974017,
    42.49935911110952,
    45.43875343077235,
    28.45211111649501,
    40.31644222904758,
    33.15774574200091,
    42.94112385639288,
    48.57542918650685,
    50.72431493096373,
    42.61042088555634,
    35.44004687602678,
    24.77106435711227
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Remove the geometry column (G) entirely: header + data.
$ws.Columns.Item(7).Delete()
